# Daily attendance processing - rotate "Recorded By" (column G) name lists
# so that the last comma-separated name moves to the front.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

# Column G holds "Recorded By". Data starts at row 2 (row 1 is the header).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Text

    if ($null -eq $value) { continue }
    if ($value -eq "") { continue }

    $parts = $value -split ",\s*"
    if ($parts.Count -gt 1) {
        $rotated = @($parts[$parts.Count - 1]) + $parts[0..($parts.Count - 2)]
        $cell.Value = [string]::Join(", ", $rotated)
    }
}
